$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7304773333333333
$ws.Range("H2").Value = 2.191432
$ws.Range("I2").Value = 0.03163269997405359
$ws.Range("J2").Value = 0.03163269997405359
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 5.469538151567999
$ws.Range("R2").Value = 49.225843364112
$ws.Range("S2").Value = 0.004379215735023403
$ws.Range("T2").Value = 0.004379215735023403
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7304773333333333
$ws.Range("H3").Value = 2.191432
$ws.Range("I3").Value = 0.03163269997405359
$ws.Range("J3").Value = 0.03163269997405359
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 23.374856833632
$ws.Range("R3").Value = 210.373711502688
$ws.Range("S3").Value = 0.01871520739286097
$ws.Range("T3").Value = 0.01871520739286097
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7304773333333333
$ws.Range("H4").Value = 2.191432
$ws.Range("I4").Value = 0.03163269997405359
$ws.Range("J4").Value = 0.03163269997405359
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 10.66410832087555
$ws.Range("R4").Value = 95.97697488787999
$ws.Range("S4").Value = 0.008538276846169217
$ws.Range("T4").Value = 0.008538276846169217
$ws.Range("I5").Value = 0.4074771110502447
$ws.Range("J5").Value = 0.4074771110502448
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 70.45593979040999
$ws.Range("R5").Value = 634.10345811369
$ws.Range("S5").Value = 0.056410934818614
$ws.Range("T5").Value = 0.05641093481861401
$ws.Range("I6").Value = 0.4074771110502447
$ws.Range("J6").Value = 0.4074771110502448
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.2410802317666319
$ws.Range("T6").Value = 0.2410802317666319
$ws.Range("I7").Value = 0.4074771110502447
$ws.Range("J7").Value = 0.4074771110502448
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 137.3698752898583
$ws.Range("R7").Value = 1236.328877608725
$ws.Range("S7").Value = 0.1099859444649988
$ws.Range("T7").Value = 0.1099859444649988
$ws.Range("G8").Value = 12.95234266666667
$ws.Range("H8").Value = 38.857028
$ws.Range("I8").Value = 0.5608901889757016
$ws.Range("J8").Value = 0.5608901889757018
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 96.98224590247199
$ws.Range("R8").Value = 872.840213122248
$ws.Range("S8").Value = 0.07764936736975865
$ws.Range("T8").Value = 0.07764936736975867
$ws.Range("G9").Value = 12.95234266666667
$ws.Range("H9").Value = 38.857028
$ws.Range("I9").Value = 0.5608901889757016
$ws.Range("J9").Value = 0.5608901889757018
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 414.467556593328
$ws.Range("R9").Value = 3730.208009339952
$ws.Range("S9").Value = 0.3318457235680622
$ws.Range("T9").Value = 0.3318457235680624
$ws.Range("G10").Value = 12.95234266666667
$ws.Range("H10").Value = 38.857028
$ws.Range("I10").Value = 0.5608901889757016
$ws.Range("J10").Value = 0.5608901889757018
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 189.0889407562244
$ws.Range("R10").Value = 1701.80046680602
$ws.Range("S10").Value = 0.1513950980378807
$ws.Range("T10").Value = 0.1513950980378807
